$wb = $excel.ActiveWorkbook

# Insert the new "dates" worksheet right after "sales" (before "dummy"),
# matching the sheet order sales / dates / dummy.
$salesSheet = $wb.Worksheets.Item("sales")
$ws = $wb.Worksheets.Add($null, $salesSheet)
$ws.Name = "dates"

# Header row text
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "sales"

# First two data rows: a literal seed value, then a relative formula
$ws.Range("A2").Value = 43101
$ws.Range("B2").Value = 1
$ws.Range("A3").Formula = "=A2+1"
$ws.Range("B3").Formula = "=B2+1"

# Remaining rows as one fill so the engine records them as shared formulas
$ws.Range("A4:A34").Formula = "=A3+1"
$ws.Range("B4:B34").Formula = "=B3+1"

# Date formatting for the date column (create this style before the bold
# header style, so the new cellXfs entries land in the same order as the
# target workbook: numFmt 15 first, bold header second)
$ws.Range("A2:A34").NumberFormat = "d-mmm-yy"

# Bold header row
$ws.Range("A1:B1").Font.Bold = $true

# Size column A to fit its (date-formatted) contents
$ws.Columns.Item(1).AutoFit() | Out-Null

# Match the page setup used by the other printable sheet in this workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# This new sheet is the active tab after the edit
$ws.Activate()
